$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.287.29"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").Value = "1.868.51"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'245.06"
$ws.Range("E5").Value = "  +4.55%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "'0.4722"
$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").Value = "'42.65"
$ws.Range("E8").Value = "  -1.30%  "

$ws.Range("D9").Value = "'0.2869"
$ws.Range("E9").Value = "  +0.19%  "

$ws.Range("D10").Value = "'0.06464"
$ws.Range("E10").Value = "  -1.60%  "

$ws.Range("D11").Value = "'20.98"
$ws.Range("E11").Value = "  -1.27%  "

$ws.Range("D12").Value = "'0.07763"
$ws.Range("E12").Value = "  -1.11%  "

$ws.Range("D13").Value = "1.878.65"
$ws.Range("E13").Value = "  +0.50%  "

$ws.Range("D14").Value = "'94.95"
$ws.Range("E14").Value = "  -1.95%  "

$ws.Range("D15").Value = "'0.7090"
$ws.Range("E15").Value = "  +2.04%  "

$ws.Range("D16").Value = "'5.092"
$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("D17").Value = "'275.49"
$ws.Range("E17").Value = "  +2.86%  "

$ws.Range("D18").Value = "30.278.63"
$ws.Range("E18").Value = "  -0.06%  "

$ws.Range("D19").Value = "'13.30"
$ws.Range("E19").Value = "  -3.77%  "

$ws.Range("D20").Value = "'0.000007540"
$ws.Range("E20").Value = "  -1.33%  "

$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("D22").Value = "2.132.04"
$ws.Range("E22").Value = "  +0.52%  "

$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "'5.208"
$ws.Range("E24").Value = "  -0.31%  "

$ws.Range("D25").Value = "'6.143"
$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("D26").Value = "'9.234"
$ws.Range("E26").Value = "  -2.00%  "

$ws.Range("D27").Value = "'165.48"
$ws.Range("E27").Value = "  -1.26%  "

$ws.Range("D28").Value = "'18.86"
$ws.Range("E28").Value = "  +0.01%  "

$ws.Range("D29").Value = "'1.900"
$ws.Range("E29").Value = "  -2.22%  "

$ws.Range("D30").Value = "'1.378"
$ws.Range("E30").Value = "  +1.27%  "

$ws.Range("D31").Value = "'0.09875"
$ws.Range("E31").Value = "  -0.46%  "

$ws.Range("E32").Value = "  +4.02%  "

$ws.Range("D33").Value = "'4.254"
$ws.Range("E33").Value = "  -2.33%  "

$ws.Range("D34").Value = "'4.020"
$ws.Range("E34").Value = "  -0.95%  "

$ws.Range("D35").Value = "'0.04763"
$ws.Range("E35").Value = "  +0.73%  "

$ws.Range("D36").Value = "'1.115"
$ws.Range("E36").Value = "  -1.75%  "

$ws.Range("D37").Value = "'0.6911"
$ws.Range("E37").Value = "  -1.39%  "

$ws.Range("D38").Value = "'2.717"
$ws.Range("E38").Value = "  +0.18%  "

$ws.Range("D39").Value = "'0.01853"
$ws.Range("E39").Value = "  -1.07%  "

$ws.Range("D40").Value = "'2.740"
$ws.Range("E40").Value = "  -2.09%  "

$ws.Range("D41").Value = "'6.291"
$ws.Range("E41").Value = "  +0.33%  "

$ws.Range("D42").Value = "'70.15"
$ws.Range("E42").Value = "  -3.58%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8410"
$ws.Range("E43").Value = "  +0.60%  "

$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'1.905"
$ws.Range("E45").Value = "  -2.62%  "

$ws.Range("D46").Value = "'0.4096"
$ws.Range("E46").Value = "  -1.81%  "

$ws.Range("D47").Value = "'101.66"
$ws.Range("E47").Value = "  -1.44%  "

$ws.Range("D48").Value = "'9.252"
$ws.Range("E48").Value = "  +1.62%  "

$ws.Range("D49").Value = "'7.061"
$ws.Range("E49").Value = "  -0.73%  "

$ws.Range("D50").Value = "'35.17"
$ws.Range("E50").Value = "  +2.04%  "

$ws.Range("D51").Value = "'917.68"
$ws.Range("E51").Value = "  -5.61%  "
